$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0. Register the "List Paragraph" style (word/styles.xml) used by the new
#    numbered Members list below.
# ---------------------------------------------------------------------------
$listStyle = $d.Styles.Add("ListParagraph", 1)
$listStyle.NameLocal = "List Paragraph"
$listStyle.BaseStyle = $d.Styles("Normal")
$listStyle.Priority = 34
$listStyle.QuickStyle = $true
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

# ---------------------------------------------------------------------------
# 1. "Group Name: DSEs" -> "Group Name: " (regular) + "DSEs" (bold)
# ---------------------------------------------------------------------------
$rngDSEs = $d.Paragraphs(2).Range.Duplicate
$rngDSEs.Find.Execute("DSEs", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngDSEs.Bold = 1

# ---------------------------------------------------------------------------
# 2. "Author: " / "Abhinav" / " " / "Bhandaram" / "." -> one run, no proofErr
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$rng3 = $d.Range($p3.Start, $p3.End)
$rng3.Text = "Author: Abhinav Bhandaram."

# ---------------------------------------------------------------------------
# 3. Turn the "Members" entries into a numbered (ListParagraph) list sharing
#    the same list (numId). ApplyNumberDefault() mints the numbering.xml
#    part + numId=1 the first time it is used; every other member paragraph
#    re-uses that numId via a direct pPr/numPr rewrite.
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.ListFormat.ApplyNumberDefault()

$p6 = $d.Paragraphs(6).Range
$xml6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Kaza</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Nikhitha</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p6.InsertXML($xml6)

$p7 = $d.Paragraphs(7).Range
$xml7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Rajashekar</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Goud </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Korakoppula</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p7.InsertXML($xml7)

# ---------------------------------------------------------------------------
# 4. New list entry "Abhinav Bhandaram." right after "Rajashekar Goud
#    Korakoppula" (still part of the same numbered list).
# ---------------------------------------------------------------------------
$p7again = $d.Paragraphs(7).Range
$collapsed = $d.Range($p7again.End, $p7again.End)
$collapsed.InsertParagraphAfter()

$p8 = $d.Paragraphs(8).Range
$xml8 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>
</w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Abhinav Bhandaram.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p8.InsertXML($xml8)

# ---------------------------------------------------------------------------
# 5. Move the "_GoBack" bookmark off the "Kaza Nikhitha" paragraph and onto
#    the trailing empty paragraph (its original home before the member list
#    grew).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(9).Range
$xml9 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p9.InsertXML($xml9)

Write-Output "done"
